$d = $word.ActiveDocument

function ReplaceAll($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# --- Simple/Traditional Chinese word-choice updates (unambiguous, global) ---
ReplaceAll "英语" "英文"
ReplaceAll " / 葡萄牙语 / 法语 / 泰语 / 越南语 / 西班牙语" " / 葡萄牙文 / 法文 / 泰文 / 越南文 / 西班牙文"
ReplaceAll "简要" "簡介"
ReplaceAll "发送给在目标国家中确认参加但尚未向我们提交文件的合作伙伴的电子邮件。 将通过 customer.io 发送" "發送給在目標國家中已回覆參加但尚未寄送文件的合作夥伴的電子郵件。 將通過 customer.io 發送"
ReplaceAll "目标受众" "目標受眾"
ReplaceAll "已邀请但尚未提交文件的合作伙伴" "尚未提交文件的受邀合作夥伴"
ReplaceAll "主题行" "主題行"
ReplaceAll "[活动名称]" "[活動名稱]"
ReplaceAll " — 您是否已提交文件？  " " — 您已提交文件了嗎？  "
ReplaceAll "不要忘记发送文件" "不要忘記傳送文檔"
ReplaceAll "[合作伙伴姓名]" "[合作夥伴姓名]"
ReplaceAll "很高兴能在即将举行的 " "很高興能在即將舉行的 "
ReplaceAll "为了确认注册，需要您在 " "為了確認註冊，需要您在 "
ReplaceAll " 之前提供以下文件：" " 之前提供以下文檔："
ReplaceAll "[插入所需文件列表]" "[插入所需文件清單]"
ReplaceAll "请将这些文件的副本发送给您的区域经理 " "請將這些文檔的副本傳送給您的區域經理 "
ReplaceAll "，邮箱地址为 " "，郵箱地址為 "
ReplaceAll "[电子邮件地址]" "[郵箱地址]"
ReplaceAll "[WHATSAPP 号码]" "[WHATSAPP 號碼]"
ReplaceAll " (WhatsApp)，以便我们做出必要的安排，包括住宿和交通。" " (WhatsApp)，以便我們做出必要的安排，包括住宿和交通。"
ReplaceAll "如有任何疑问，请联系您的区域经理。" "如有任何疑問，請聯繫您的區域經理。"
ReplaceAll "期待在那里见到您！" "期待在那裡見到您！"
ReplaceAll "为了确保您在此活动中获得最佳体验，我们需要您在 " "為了確保您在此次活動中擁有最佳體驗，我們需要您在 "
ReplaceAll "请回复此电子邮件，附上这些文件的副本，以便我们能为您做出必要的安排，包括住宿和交通。" "請回覆此電子郵件，附上這些文檔的副本，以便我們做出必要的安排，包括住宿和交通。"
ReplaceAll "实时聊天" "即時聊天"
ReplaceAll "如果您有任何问题，请联系您的区域经理 " "如有任何疑問，請聯繫您的區域經理 "

# --- Ambiguous spots: same source text resolves to different targets depending on
#     which paragraph it is in, so these are scoped to a single paragraph each. ---

# "简要：" -> "簡介:" (colon loses the following space)
$d.Paragraphs(5).Range.Find.Execute("：", $true, $false, $false, $false, $false, $true, 1, $false, ":", 2) | Out-Null

# "目标受众：" -> "目標受眾:" (colon, no space)
$d.Paragraphs(8).Range.Find.Execute("：", $true, $false, $false, $false, $false, $true, 1, $false, ":", 2) | Out-Null

# "主题行：" -> "主題行: " (colon + trailing space) -- subject line 1
$d.Paragraphs(12).Range.Find.Execute("：", $true, $false, $false, $false, $false, $true, 1, $false, ": ", 2) | Out-Null

# "主题行：" -> "主題行: " (colon + trailing space) -- subject line 2
$d.Paragraphs(29).Range.Find.Execute("：", $true, $false, $false, $false, $false, $true, 1, $false, ": ", 2) | Out-Null

# " 见到您。 " stays the same in paragraph 17 (first block) -- no textual change needed there
# beyond what global replacements already covered.

# " 见到您。 ‘" -> " 見到您。 ‘" in paragraph 35 (second block, trailing curly quote kept)
$d.Paragraphs(35).Range.Find.Execute(" 见到您。 ‘", $true, $false, $false, $false, $false, $true, 1, $false, " 見到您。 ‘", 2) | Out-Null
$d.Paragraphs(17).Range.Find.Execute(" 见到您。 ", $true, $false, $false, $false, $false, $true, 1, $false, " 見到您。 ", 2) | Out-Null

# --- Comment text update (Traditional Chinese) ---
$d.Comments(1).Range.Text = "選擇任一"
